$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells are stored as plain text in the workbook ("Price" column uses
# dotted/grouped formatting like "60.626.52" or fixed-precision strings like "0.998").
# For values that otherwise look like a valid number, prefix with a literal single
# quote (Excel's text-entry marker) so COM stores them as text instead of coercing
# them to a Double (which would also strip trailing zeros, e.g. "18.50" -> 18.5).

$ws.Range('D2').Value = '60.626.52'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '2.343.69'
$ws.Range('E3').Value = '  -2.88%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '''544.06'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '''135.94'
$ws.Range('E6').Value = '  -5.38%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  -8.83%  '
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').Value = '''24.48'
$ws.Range('E14').Value = '  -3.96%  '
$ws.Range('D15').Value = '2.767.00'
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').Value = '60.531.15'
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('E17').Value = '  -2.36%  '
$ws.Range('D18').Value = '2.343.46'
$ws.Range('E18').Value = '  -2.73%  '
$ws.Range('D19').Value = '''10.56'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').Value = '''317.67'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').Value = '''6.54'
$ws.Range('E22').Value = '  -3.54%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  -4.29%  '
$ws.Range('E26').Value = '  +7.88%  '
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '''497.69'
$ws.Range('E28').Value = '  -2.68%  '
$ws.Range('E29').Value = '  -4.78%  '
$ws.Range('D30').Value = '0.0₃0857'
$ws.Range('E30').Value = '  -9.08%  '
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('E32').Value = '  -3.49%  '
$ws.Range('E33').Value = '  -3.58%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').Value = '''0.375'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '''18.50'
$ws.Range('E37').Value = '  +2.53%  '
$ws.Range('E38').Value = '  -5.62%  '
$ws.Range('D39').Value = '''1.80'
$ws.Range('E39').Value = '  +5.16%  '
$ws.Range('D40').Value = '''141.14'
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').Value = '''141.46'
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').Value = '''2.06'
$ws.Range('E45').Value = '  -5.69%  '
$ws.Range('D46').Value = '''0.0510'
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').Value = '''18.94'
$ws.Range('E47').Value = '  -7.65%  '
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('E51').Value = '  -2.33%  '
